$wb = $excel.ActiveWorkbook

# --- Insert new "DisableEffects" sheet right after "Global" ---
$global = $wb.Worksheets.Item("Global")
$disableEffects = $wb.Worksheets.Add($null, $global)
$disableEffects.Name = "DisableEffects"

$disableEffects.Range("A1").Value = "KEY"
$disableEffects.Range("B1").Value = "ENGLISH"
$disableEffects.Range("C1").Value = "KOREAN"
$disableEffects.Range("D1").Value = "SPANISH"

$disableEffects.Range("A2").Value = "NAME"
$disableEffects.Range("B2").Value = "Disable Effects"

$disableEffects.Range("A3").Value = "DESCRIPTION"
$disableEffects.Range("B3").Value = "Disables certain effects to improve performance."

$disableEffects.Range("A4").Value = "FILTER"
$disableEffects.Range("B4").Value = "Disable VFX filters (Grayscale, Arcade, etc.)"

$disableEffects.Range("A5").Value = "BLOOM"
$disableEffects.Range("B5").Value = "Disable bloom"

$disableEffects.Range("A6").Value = "FLASH"
$disableEffects.Range("B6").Value = "Disable screen flashes"

$disableEffects.Range("A7").Value = "HALL_OF_MIRRORS"
$disableEffects.Range("B7").Value = 'Disable "Hall of Mirrors" effect'

$disableEffects.Range("A8").Value = "SCREEN_SHAKE"
$disableEffects.Range("B8").Value = "Disable screen shake"

# --- Add missing Korean translations to "KeyLimiter" sheet ---
$keyLimiter = $wb.Worksheets.Item("KeyLimiter")
$keyLimiter.Range("C8").Value = "CLS (커스텀 레벨 선택) 에서 키 제한하기"
$keyLimiter.Range("C9").Value = "메인 메뉴에서 키 제한하기"

# --- Restore original active sheet/selection (Global was active before the edit) ---
$global.Activate()
[void]$global.Range("A1").Select()
